$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# The blank/empty treatment count values in column X (rows 7-14) are set to -1.
for ($row = 7; $row -le 14; $row++) {
    $ws.Cells.Item($row, 24).Value = -1
}

# Reset the view: select A1 so the sheet no longer shows the scrolled-to
# "D1" top-left position with "W2" selected.
$ws.Activate()
$ws.Range("A1").Select()
